$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.293.05"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.383.28"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.98"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.22"
$ws.Range("E6").Value = "  -6.33%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.383.79"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.46"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "3.960.91"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "3.394.27"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").Value = "60.413.74"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.96"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -6.52%  "
$ws.Range("D27").Value = "3.528.81"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -5.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("E33").Value = "  -7.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.60"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").Value = "3.413.95"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "167.61"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("E40").Value = "  -7.22%  "
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.76"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.43"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").Value = "2.514.78"
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("E51").Value = "  -3.87%  "
